$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Kit-components table (Table 3): normalize whitespace handling on the
#    existing cell text runs by doing a full-document Find/ReplaceAll of
#    each distinct string with itself. This rewrites the backing <w:t>
#    nodes so Word no longer stamps them with xml:space="preserve" (none of
#    these values have leading/trailing whitespace).
# ---------------------------------------------------------------------------
$wdReplaceAll = 2
$wdFindContinue = 1

function Normalize-Text($text) {
    $rng = $d.Content
    [void]$rng.Find.Execute($text, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $text, $wdReplaceAll)
}

$kitTexts = @(
    "Anti-Mouse Klk1 Pre-coated 96-well Strip Microplate",
    "12 strips of 8 wells",
    "Return unused wells to the foil pouch. Reseal along the entire edge of the zip-seal. May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit.",
    "Mouse Klk1 Standard",
    "10 ng/tube",
    "Discard the Klk1 stock solution after 12 hours at 4°C. May be stored at -20°C for 48 hours.",
    "Mouse Klk1 Biotinylated Antibody (100x)",
    "100 µl",
    "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit.",
    "Avidin-Biotin-Peroxidase Complex (100x)",
    "Sample Diluent",
    "30 ml",
    "Antibody Diluent",
    "12 ml",
    "Avidin-Biotin-Peroxidase Diluent",
    "1",
    "2"
)

foreach ($txt in $kitTexts) {
    Normalize-Text $txt
}

# ---------------------------------------------------------------------------
# 2) Append four new rows to the kit-components table (Table 3): Color
#    Developing Reagent (TMB), Stop Solution, Wash Buffer (25x), Plate
#    Sealers (last cell of the last row intentionally left blank).
# ---------------------------------------------------------------------------
$kitTable = $d.Tables.Item(3)

$newKitRows = @(
    @("Color Developing Reagent (TMB)", "1", "10 ml", "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."),
    @("Stop Solution", "1", "10 ml", "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."),
    @("Wash Buffer (25x)", "1", "20 ml", "May be stored for up to 1 month at 4°C provided this is within the expiration date of the kit."),
    @("Plate Sealers", "4", "Piece", $null)
)

foreach ($rowValues in $newKitRows) {
    $row = $kitTable.Rows.Add()
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $val = $rowValues[$c]
        if ($null -ne $val) {
            $row.Cells.Item($c + 1).Range.Text = $val
        }
    }
}

# ---------------------------------------------------------------------------
# 3) Intra-assay precision table (Table 5): update the measured values for
#    the three sample rows.
# ---------------------------------------------------------------------------
$intraTable = $d.Tables.Item(5)

$intraValues = @(
    @("24", "145", "10.15", "7.0%"),
    @("24", "618", "49.44", "8.0%"),
    @("24", "1426", "128.34", "9.0%")
)

for ($r = 0; $r -lt $intraValues.Length; $r++) {
    $rowVals = $intraValues[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $cell = $intraTable.Cell($r + 2, $c + 2)
        $cell.Range.Text = $rowVals[$c]
    }
}

# ---------------------------------------------------------------------------
# 4) Remove the centered paragraph alignment from the data rows (not the
#    header row) of both the intra-assay (Table 5) and inter-assay
#    (Table 6) precision tables. Setting Alignment back to the default
#    (wdAlignParagraphLeft = 0) drops the now-redundant <w:pPr><w:jc/></w:pPr>
#    wrapper that the template had baked in.
# ---------------------------------------------------------------------------
$wdAlignParagraphLeft = 0

foreach ($tableIndex in 5, 6) {
    $table = $d.Tables.Item($tableIndex)
    for ($r = 2; $r -le $table.Rows.Count; $r++) {
        for ($c = 1; $c -le $table.Columns.Count; $c++) {
            $cell = $table.Cell($r, $c)
            foreach ($p in $cell.Range.Paragraphs) {
                $p.Alignment = $wdAlignParagraphLeft
            }
        }
    }
}

Write-Output "edit complete"
